$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
# Columns B, E, F, H and I no longer hold data in the new layout, so clear
# them (their former shared-string values are being removed from the
# workbook entirely).
$ws.Range("B2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()

# C2 now simply holds the file extension, and D2 becomes a formula that
# concatenates the path (A2), an (empty) document name (B2) and the
# extension (C2).
$ws.Range("C2").Value = ".docx"
$ws.Range("D2").Formula = "=A2&B2&C2"

# G2 keeps its date-formatted style but no longer has a value.
$ws.Range("G2").Value = $null

# --- Row 3 (new) -----------------------------------------------------------
# Mirrors row 2: CAMINHO in A3, extension in C3, concatenation formula in
# D3, and G3 reuses G2's (date) number format without a value.
$ws.Range("A3").Value = "modelosPadrao\"
$ws.Range("C3").Value = ".docx"
$ws.Range("D3").Formula = "=A3&B3&C3"

$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)  # xlPasteFormats - reuse G2's style without touching its value
$excel.CutCopyMode = $false

# Column G widened to fit the new (wider) content.
$ws.Columns.Item(7).ColumnWidth = 17.5

# The active selection moved to D1.
$ws.Range("D1").Select()
